$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: column A switches from date-only format back to the full datetime
# format (style goes from the "YYYY-MM-DD" style to the "YYYY-MM-DD HH:MM:SS" style).
$ws.Range("A11").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 12: brand-new row of data, mirroring the structure of the existing rows.
$ws.Range("A12").Value = 45856
$ws.Range("A12").NumberFormat = "YYYY-MM-DD"

$ws.Range("B12").Value = "agora"
$ws.Range("C12").Value = "agora"
$ws.Range("D12").Value = "2025-07-18 14:07:57"
$ws.Range("E12").Value = "2025-07-18 14:07:57"
$ws.Range("F12").Value = "2025-07-18 14:07:59"
$ws.Range("G12").Value = "2025-07-18 14:08:00"
$ws.Range("H12").Value = "2025-07-18 14:08:01"
$ws.Range("I12").Value = "2025-07-18 14:08:02"
$ws.Range("J12").Value = "2025-07-18 14:08:03"
$ws.Range("K12").Value = "0:00:01"
$ws.Range("L12").Value = "0:00:00"
$ws.Range("M12").Value = "0:00:06"
$ws.Range("N12").Value = ""
$ws.Range("O12").Value = "2025-07-18 14:08:05"
$ws.Range("P12").Value = "2025-07-18 14:08:05"
$ws.Range("Q12").Value = "2025-07-18 14:08:07"
$ws.Range("R12").Value = "2025-07-18 14:08:07"
$ws.Range("S12").Value = "0:00:02"
$ws.Range("T12").Value = "0:00:01"
$ws.Range("U12").Value = "0:00:03"
$ws.Range("V12").Value = "0:00:01"
$ws.Range("W12").Value = "2025-07-18 14:08:04"
